$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Title ---
Replace-Text "Delving into the Enigma of Dreams" "The Game of Numbers: Exploring the Fascinating World of Mathematics"

# --- Author name ---
Replace-Text " Evelyn Carter" " Anna Jackson"

# --- Email / handle ---
Replace-Text "evelyn" "anna"
Replace-Text "carter@arcaneresearchuni" "jackson.88@academics.ville"

# --- Body paragraph 1 (sz 24), first block before the first pair of breaks ---
Replace-Text "Throughout the ages, dreams have captivated and perplexed humankind, transcending cultural, geographical, and temporal boundaries" "In a world enveloped by complexities, Mathematics remains a luminous beacon of universal order"

Replace-Text " Emerging from the enigmatic depths of the subconscious, dreams propel us into a realm where reality and imagination intertwine" " It is the language of patterns, the murmur of numbers harmonizing our existence, and the exquisite dance of reason"

Replace-Text " As fleeting as morning mist, they vanish upon waking, leaving behind fragments of memories that linger in our consciousness" " From the smallest atom to the celestial vastness, Mathematics serves as the pervasive tome of nature"

Replace-Text " Yet, despite their ephemeral nature, dreams persist as an enduring enigma, beckoning us to unravel their mysteries" " Its principles anchor our everyday lives, weaving a delicate tapestry of logic, measurement, and abstraction. This journey into the realm of Mathematics will beckon us to comprehend its extraordinary essence"

# --- Body paragraph 1, second block (after first pair of breaks) ---
Replace-Text "These nocturnal journeys stir curiosities from diverse disciplines" "History has borne witness to the profound impact of Mathematics in shaping civilizations"

Replace-Text " Neuroscientists tirelessly probe the mechanisms that orchestrate the intricate narrative of dreams, unraveling the complex interplay of brain regions that craft these ethereal landscapes" " From the intricate writings of ancient Babylonian mathematicians to the transformative discoveries of Greek prodigies like Pythagoras and Euclid, Mathematics has held an unyielding allure for intellectual seekers"

Replace-Text " Psychologists delve into the psyche's hidden recesses, deciphering the symbolism and messages encoded within dream imagery, revealing glimpses of our inner selves" " Through the ages, luminaries like Archimedes, Newton, and Einstein pushed the boundaries of mathematical knowledge, unveiling secrets of the cosmos. This progression of thought, a seamless dialogue between brilliant minds, elucidates the enduring quest for unlocking the enigmatic mysteries of existence"

# --- Body paragraph 1, third block (after second pair of breaks) ---
Replace-Text "Across civilizations, indigenous cultures perceive dreams as portals to the supernatural, imparting divine messages or glimpses into parallel realms" "In the pursuit of knowledge, Mathematics offers not just rigorous methods but also a unique kind of wisdom"

Replace-Text " Artistic souls find inspiration in the surreal tapestries of dreams, weaving their enigmatic tales into masterpieces" " It instills discipline, sharpens logical reasoning, and encourages problem-solving prowess"

Replace-Text " Poets and writers craft verses and narratives imbued with dream-like imagery, capturing the ethereal essence of these elusive experiences" " Beyond its functional utility, Mathematics unravels a breathtaking saga of patterns, symmetries, and interconnections--a true symphony of ideas. Its aesthetics, elegance, and universality stand as testaments to the inherent beauty and order underpinning the universe. Engage with mathematics is to embark on a profound intellectual odyssey, an exploration of a realm as vast and infinite as the cosmos itself"

# --- Summary paragraph ---
Replace-Text "Dreams, enigmatic voyages of the unconscious mind, hold a mirror to our inner selves, revealing hidden corners of our psyche" "Mathematics encompasses the study of patterns and relationships, forming the bedrock of our logical and conceptual understanding of the world"

Replace-Text " Their multidisciplinary significance spans neuroscience, psychology, culture, arts, and more" " It permeates our lives, from commerce to science, with its principles applied in diverse fields"

Replace-Text " As we continue to explore the realms of dreams, we uncover their capacity to illuminate our conscious lives, providing insights into our emotions, motivations, and subconscious desires" " Mathematics is a multifaceted discipline encompassing logic, measurement, number theory, geometry, probability, and statistics"

Replace-Text " The study of dreams offers an avenue for personal growth, creative inspiration, and a profound understanding of the intricate workings of the human mind" " Its history is storied, featuring brilliant minds who expanded our mathematical horizons, laying the foundation for further exploration. Furthermore, Mathematics cultivates invaluable cognitive skills, fostering critical thinking, problem-solving abilities, and a precise grasp of concepts. Beyond practicality, it presents an artistic, almost poetic, side, revealing the wonder and beauty inherent in abstract thought"

# --- Add trailing empty paragraph after the Summary paragraph ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
